# Auto-generated: apply cryptos list update (commit: "Updated cryptos list on Sun Jul 14 03:20:12 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.478.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.72%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.177.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.60%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.60%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.23%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.30%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.113"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.54%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.430"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.724.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.56%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.140"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "

# Row 15
$ws.Range("E15").Value = "  +3.58%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.513.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.182.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.75%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.92%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.12%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.529"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.10%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.61%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.24%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +16.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0895"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.75%  "

# Row 29
$ws.Range("E29").Value = "  +1.47%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.53%  "

# Row 31
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.65%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.63%  "

# Row 33
$ws.Range("E33").Value = "  -0.60%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.12%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.53%  "

# Row 36
$ws.Range("E36").Value = "  +3.82%  "

# Row 37
$ws.Range("E37").Value = "  +6.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.729.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.95%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.29%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.56%  "

# Row 42
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.724"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.19%  "

# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.89%  "

# Row 44
$ws.Range("E44").Value = "  +7.40%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.220.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.60%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.48%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.100"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.16%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.14%  "

# Row 50
$ws.Range("E50").Value = "  -0.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.768"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.31%  "

